$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $trimmed = $d.Range($r.Start, $r.End - 1)
    $trimmed.Text = $newText
}

# QUOTAS.cod is column 5, QUOTAS.hake is column 6.
# Row 1 is the header; rows 2-15 correspond to BE, DK, DE, EE, IE, ES, FR, LV, LT, NL, PL, PT, FI, SE.

Set-CellText $t 2 5 "0.487"   # BE QUOTAS.cod  0.488 -> 0.487
Set-CellText $t 2 6 "0.492"   # BE QUOTAS.hake 0.494 -> 0.492

Set-CellText $t 3 5 "0.174"   # DK QUOTAS.cod  0.175 -> 0.174
Set-CellText $t 3 6 "0.528"   # DK QUOTAS.hake 0.529 -> 0.528

Set-CellText $t 4 5 "0.520"   # DE QUOTAS.cod  0.521 -> 0.520
Set-CellText $t 4 6 "0.489"   # DE QUOTAS.hake 0.490 -> 0.489

Set-CellText $t 5 5 "0.479"   # EE QUOTAS.cod  0.480 -> 0.479
Set-CellText $t 5 6 "0.000"   # EE QUOTAS.hake 0.489 -> 0.000

Set-CellText $t 6 5 "0.459"   # IE QUOTAS.cod  0.460 -> 0.459
Set-CellText $t 6 6 "0.000"   # IE QUOTAS.hake 0.515 -> 0.000

Set-CellText $t 7 5 "0.638"   # ES QUOTAS.cod  0.639 -> 0.638

Set-CellText $t 8 5 "0.475"   # FR QUOTAS.cod  0.476 -> 0.475

Set-CellText $t 9 5 "0.457"   # LV QUOTAS.cod  0.458 -> 0.457
Set-CellText $t 9 6 "0.000"   # LV QUOTAS.hake 0.489 -> 0.000

Set-CellText $t 10 5 "0.466"  # LT QUOTAS.cod  0.467 -> 0.466
Set-CellText $t 10 6 "0.000"  # LT QUOTAS.hake 0.489 -> 0.000

Set-CellText $t 11 5 "0.484"  # NL QUOTAS.cod  0.485 -> 0.484
Set-CellText $t 11 6 "0.492"  # NL QUOTAS.hake 0.493 -> 0.492

Set-CellText $t 12 6 "0.000"  # PL QUOTAS.hake 0.489 -> 0.000

Set-CellText $t 13 5 "0.519"  # PT QUOTAS.cod  0.520 -> 0.519
Set-CellText $t 13 6 "0.532"  # PT QUOTAS.hake 0.533 -> 0.532

Set-CellText $t 14 5 "0.459"  # FI QUOTAS.cod  0.460 -> 0.459
Set-CellText $t 14 6 "0.000"  # FI QUOTAS.hake 0.489 -> 0.000

Set-CellText $t 15 5 "0.364"  # SE QUOTAS.cod  0.365 -> 0.364
Set-CellText $t 15 6 "0.489"  # SE QUOTAS.hake 0.490 -> 0.489
